$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with new columns D and E, matching the style of C1 (bold, bordered, centered)
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4

# Bulk-write updated matrix values for rows 2-67, columns B:E
$data = New-Object 'object[,]' 66,4
$data[0,0] = -0.323716723852815
$data[0,1] = -0.3093107257104907
$data[0,2] = -0.2952607650692658
$data[0,3] = -0.2819917716829335
$data[1,0] = 0.198459336637776
$data[1,1] = 0.2076865250416704
$data[1,2] = 0.2176271504216085
$data[1,3] = 0.2275052025583008
$data[2,0] = 0.1308398275189065
$data[2,1] = 0.1425296526411839
$data[2,2] = 0.1551052188821808
$data[2,3] = 0.1675804533426762
$data[3,0] = -0.08598760563481117
$data[3,1] = -0.07633509935382571
$data[3,2] = -0.06723701352770255
$data[3,3] = -0.05904631302620567
$data[4,0] = 0.1257264185200478
$data[4,1] = 0.1418282182027905
$data[4,2] = 0.157312441416569
$data[4,3] = 0.171410369889013
$data[5,0] = -0.462794068214327
$data[5,1] = -0.4528029662052117
$data[5,2] = -0.4425278526256859
$data[5,3] = -0.4323507975944151
$data[6,0] = -0.2722784098748927
$data[6,1] = -0.260071556078068
$data[6,2] = -0.2472640206060024
$data[6,3] = -0.2345191521032665
$data[7,0] = -0.419160127918516
$data[7,1] = -0.4090605937397951
$data[7,2] = -0.3970385238284929
$data[7,3] = -0.3839031396493626
$data[8,0] = 0.3437353167829754
$data[8,1] = 0.358506294481997
$data[8,2] = 0.3723048670123302
$data[8,3] = 0.3846022088222938
$data[9,0] = -0.2495262058236078
$data[9,1] = -0.2401966792289666
$data[9,2] = -0.2309574200385598
$data[9,3] = -0.2222111348588156
$data[10,0] = -0.1090023623403518
$data[10,1] = -0.08648358856128772
$data[10,2] = -0.06639254954047909
$data[10,3] = -0.04894333054410106
$data[11,0] = -0.01412646471737876
$data[11,1] = -0.0135888336771998
$data[11,2] = -0.0126784057604223
$data[11,3] = -0.0118948220840237
$data[12,0] = 0.01271392789801482
$data[12,1] = 0.05244956046941167
$data[12,2] = 0.08869275783823745
$data[12,3] = 0.1205847614118937
$data[13,0] = -0.0570807994248989
$data[13,1] = -0.01300978491836207
$data[13,2] = 0.02601119802271654
$data[13,3] = 0.05933515250160003
$data[14,0] = 0.2909665680821449
$data[14,1] = 0.3511649265855035
$data[14,2] = 0.4050524548647043
$data[14,3] = 0.4518845707694348
$data[15,0] = 0.5685316785691956
$data[15,1] = 0.6020326650123404
$data[15,2] = 0.6307855370710789
$data[15,3] = 0.6540451976054349
$data[16,0] = 0.0514454229502002
$data[16,1] = 0.03866759896343645
$data[16,2] = 0.02762113004102702
$data[16,3] = 0.01765353438448207
$data[17,0] = 0.3784449178779835
$data[17,1] = 0.3942978762976164
$data[17,2] = 0.4085839877825673
$data[17,3] = 0.4203445844300895
$data[18,0] = 0.2404340588128443
$data[18,1] = 0.2970493814345418
$data[18,2] = 0.347113586368768
$data[18,3] = 0.3898407705456988
$data[19,0] = 0.4641284581422794
$data[19,1] = 0.5205383953746511
$data[19,2] = 0.56878821272325
$data[19,3] = 0.6089081245250895
$data[20,0] = 0.3442710520518817
$data[20,1] = 0.3797409485700736
$data[20,2] = 0.410459550677119
$data[20,3] = 0.4357712697100385
$data[21,0] = -0.09234164219002516
$data[21,1] = -0.06102057954108653
$data[21,2] = -0.03429578103519461
$data[21,3] = -0.01251349659053289
$data[22,0] = 4.577738315524726
$data[22,1] = 4.600142545217028
$data[22,2] = 4.543939092626456
$data[22,3] = 4.425746139532285
$data[23,0] = 0.5201919331051155
$data[23,1] = 0.4798567206283482
$data[23,2] = 0.4490839207864836
$data[23,3] = 0.4213661139698385
$data[24,0] = 0.3840934113596692
$data[24,1] = 0.3525363144340392
$data[24,2] = 0.3324494749365993
$data[24,3] = 0.3122743053466845
$data[25,0] = 0.3506188273330846
$data[25,1] = 0.310901809741786
$data[25,2] = 0.2802978279369414
$data[25,3] = 0.250300738809979
$data[26,0] = 1.107158495914431
$data[26,1] = 1.066890520460006
$data[26,2] = 1.034384186030383
$data[26,3] = 1.003125729348551
$data[27,0] = 5.891534824325058
$data[27,1] = 5.468790468361631
$data[27,2] = 5.054774525521304
$data[27,3] = 4.654705059835492
$data[28,0] = 1.02554562774442
$data[28,1] = 0.9649696633972918
$data[28,2] = 0.918192690939621
$data[28,3] = 0.8754001097957069
$data[29,0] = -0.1632159372585917
$data[29,1] = -0.2243439863291791
$data[29,2] = -0.2732557685607494
$data[29,3] = -0.3160735247271406
$data[30,0] = 0.779564758101126
$data[30,1] = 0.7474304657757006
$data[30,2] = 0.7241680675817791
$data[30,3] = 0.7011752861778168
$data[31,0] = 0.9300366235590389
$data[31,1] = 0.8993922506230235
$data[31,2] = 0.8769816692917198
$data[31,3] = 0.8559540949582878
$data[32,0] = -0.6329906594375985
$data[32,1] = -0.6672725585690531
$data[32,2] = -0.6938528429753864
$data[32,3] = -0.7189302401520412
$data[33,0] = 0.8222461609779662
$data[33,1] = 0.8131609638990331
$data[33,2] = 0.8053063940955221
$data[33,3] = 0.7980393932961181
$data[34,0] = 0.7849608595243773
$data[34,1] = 0.767376525621137
$data[34,2] = 0.7524983959128847
$data[34,3] = 0.7395860261111458
$data[35,0] = 0.7610204137432285
$data[35,1] = 0.7413393640081756
$data[35,2] = 0.7244343730520797
$data[35,3] = 0.7097215971725575
$data[36,0] = 0.735657008090314
$data[36,1] = 0.7166327809518839
$data[36,2] = 0.7001331228593332
$data[36,3] = 0.6853823668620735
$data[37,0] = 0.584655617669402
$data[37,1] = 0.5805097688759314
$data[37,2] = 0.577357603485031
$data[37,3] = 0.5746130955545922
$data[38,0] = 0.7535855833296246
$data[38,1] = 0.7508300768615521
$data[38,2] = 0.7484278463438755
$data[38,3] = 0.7457993942262386
$data[39,0] = 0.5628614257754467
$data[39,1] = 0.5545927399189651
$data[39,2] = 0.5486144725315969
$data[39,3] = 0.5440646841984613
$data[40,0] = 0.7162239923506889
$data[40,1] = 0.6890526931284596
$data[40,2] = 0.6651710945092365
$data[40,3] = 0.6439435367320129
$data[41,0] = 0.7235113542228019
$data[41,1] = 0.708607646620943
$data[41,2] = 0.6960804399275728
$data[41,3] = 0.6851437963668882
$data[42,0] = 0.6818361315274433
$data[42,1] = 0.6736650154960108
$data[42,2] = 0.667378594007934
$data[42,3] = 0.6622459522214693
$data[43,0] = 0.6758304549915243
$data[43,1] = 0.6573158281272644
$data[43,2] = 0.641827598554511
$data[43,3] = 0.6288028937328897
$data[44,0] = -1.262668264875522
$data[44,1] = -1.265986458938779
$data[44,2] = -1.269030954015268
$data[44,3] = -1.271612010106198
$data[45,0] = -0.9772038869607644
$data[45,1] = -0.9814853793243534
$data[45,2] = -0.9853535474519329
$data[45,3] = -0.9886784588588577
$data[46,0] = -0.8701089387968857
$data[46,1] = -0.8755392624003318
$data[46,2] = -0.8796802197049783
$data[46,3] = -0.8826439959900664
$data[47,0] = -0.6410403872334991
$data[47,1] = -0.6439722445150895
$data[47,2] = -0.6459379106875854
$data[47,3] = -0.6471703845323308
$data[48,0] = -0.04929204634115188
$data[48,1] = -0.05022248240209715
$data[48,2] = -0.05079377333725028
$data[48,3] = -0.05141931640286551
$data[49,0] = -0.8617558046222666
$data[49,1] = -0.8664839627003759
$data[49,2] = -0.8702238894516947
$data[49,3] = -0.8729639076070035
$data[50,0] = -0.8617558046222666
$data[50,1] = -0.8664839627003759
$data[50,2] = -0.8702238894516947
$data[50,3] = -0.8729639076070035
$data[51,0] = -1.087251285867312
$data[51,1] = -1.099638969396312
$data[51,2] = -1.109998231643979
$data[51,3] = -1.118399283333831
$data[52,0] = -0.1888127969028661
$data[52,1] = -0.1879309031566863
$data[52,2] = -0.1866147396056724
$data[52,3] = -0.185234103052067
$data[53,0] = -0.9985492423653456
$data[53,1] = -1.001920247713072
$data[53,2] = -1.004967698858919
$data[53,3] = -1.007640524675248
$data[54,0] = -0.8909248614359153
$data[54,1] = -0.8865783913471785
$data[54,2] = -0.8834885479073487
$data[54,3] = -0.8814316125723315
$data[55,0] = -0.9483704188240569
$data[55,1] = -0.9350879674765653
$data[55,2] = -0.9233701041371302
$data[55,3] = -0.9133006486670864
$data[56,0] = -1.152693352750958
$data[56,1] = -1.129367752215911
$data[56,2] = -1.109053444590487
$data[56,3] = -1.091526292210452
$data[57,0] = -0.860144294385646
$data[57,1] = -0.8481378087766493
$data[57,2] = -0.8369545860399286
$data[57,3] = -0.8266710382808649
$data[58,0] = -0.5178204024029087
$data[58,1] = -0.4975904977737807
$data[58,2] = -0.4792345968531878
$data[58,3] = -0.4629652518379018
$data[59,0] = 0.3687250643742834
$data[59,1] = 0.370916657261192
$data[59,2] = 0.3735755437270625
$data[59,3] = 0.3760671398139814
$data[60,0] = -1.232938856462481
$data[60,1] = -1.212315704557025
$data[60,2] = -1.194278238216717
$data[60,3] = -1.178811243742369
$data[61,0] = -0.76348863415763
$data[61,1] = -0.7377023705146881
$data[61,2] = -0.7135873566907546
$data[61,3] = -0.6910750878967837
$data[62,0] = -0.9092640817244033
$data[62,1] = -0.9026255638979531
$data[62,2] = -0.8951688976081931
$data[62,3] = -0.8872346792891966
$data[63,0] = -0.1304048256125618
$data[63,1] = -0.1111563709664542
$data[63,2] = -0.09375412830982176
$data[63,3] = -0.07847335130381627
$data[64,0] = -0.8170160309513503
$data[64,1] = -0.7926290785388965
$data[64,2] = -0.7719040014687291
$data[64,3] = -0.7547517506085408
$data[65,0] = -0.7972691155758459
$data[65,1] = -0.7645593498930221
$data[65,2] = -0.7375200940845188
$data[65,3] = -0.7158653422760036
$ws.Range("B2:E67").Value = $data
